$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "self" sending/target cluster pairs (bottom-to-top to keep indices stable)
$ws.Rows.Item(17).EntireRow.Delete()
$ws.Rows.Item(12).EntireRow.Delete()
$ws.Rows.Item(7).EntireRow.Delete()
$ws.Rows.Item(2).EntireRow.Delete()

# Refresh the TPM-derived values for the remaining 12 rows (recomputed stats)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Trf"
$ws.Range("C2").Value = "Tfr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6220463333333334
$ws.Range("H2").Value = 1.866139
$ws.Range("I2").Value = 0.00505260120118785
$ws.Range("J2").Value = 0.00505260120118785
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.771067666666666
$ws.Range("N2").Value = 5.313203
$ws.Range("O2").Value = 0.9584282089297276
$ws.Range("P2").Value = 0.9584282089297276
$ws.Range("Q2").Value = 1.101686148135222
$ws.Range("R2").Value = 9.915175333216999
$ws.Range("S2").Value = 0.004842555519690661
$ws.Range("T2").Value = 0.004842555519690661

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Trf"
$ws.Range("C3").Value = "Tfr2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6220463333333334
$ws.Range("H3").Value = 1.866139
$ws.Range("I3").Value = 0.00505260120118785
$ws.Range("J3").Value = 0.00505260120118785
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05687433333333333
$ws.Range("N3").Value = 0.170623
$ws.Range("O3").Value = 0.03077802528761218
$ws.Range("P3").Value = 0.03077802528761218
$ws.Range("Q3").Value = 0.03537847051077778
$ws.Range("R3").Value = 0.318406234597
$ws.Range("S3").Value = 0.0001555090875383793
$ws.Range("T3").Value = 0.0001555090875383793

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Trf"
$ws.Range("C4").Value = "Tfr2"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6220463333333334
$ws.Range("H4").Value = 1.866139
$ws.Range("I4").Value = 0.00505260120118785
$ws.Range("J4").Value = 0.00505260120118785
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01994566666666667
$ws.Range("N4").Value = 0.059837
$ws.Range("O4").Value = 0.01079376578266031
$ws.Range("P4").Value = 0.01079376578266031
$ws.Range("Q4").Value = 0.01240712881588889
$ws.Range("R4").Value = 0.111664159343
$ws.Range("S4").Value = [double]"5.45365939588098e-05"
$ws.Range("T4").Value = [double]"5.45365939588098e-05"

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Trf"
$ws.Range("C5").Value = "Tfr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.60581866666667
$ws.Range("H5").Value = 58.81745600000001
$ws.Range("I5").Value = 0.1592492032139157
$ws.Range("J5").Value = 0.1592492032139158
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.771067666666666
$ws.Range("N5").Value = 5.313203
$ws.Range("O5").Value = 0.9584282089297276
$ws.Range("P5").Value = 0.9584282089297276
$ws.Range("Q5").Value = 34.72323151906311
$ws.Range("R5").Value = 312.509083671568
$ws.Range("S5").Value = 0.1526289286097995
$ws.Range("T5").Value = 0.1526289286097995

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Trf"
$ws.Range("C6").Value = "Tfr2"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 19.60581866666667
$ws.Range("H6").Value = 58.81745600000001
$ws.Range("I6").Value = 0.1592492032139157
$ws.Range("J6").Value = 0.1592492032139158
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05687433333333333
$ws.Range("N6").Value = 0.170623
$ws.Range("O6").Value = 0.03077802528761218
$ws.Range("P6").Value = 0.03077802528761218
$ws.Range("Q6").Value = 1.115067866120889
$ws.Range("R6").Value = 10.035610795088
$ws.Range("S6").Value = 0.00490137600354999
$ws.Range("T6").Value = 0.00490137600354999

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Trf"
$ws.Range("C7").Value = "Tfr2"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 19.60581866666667
$ws.Range("H7").Value = 58.81745600000001
$ws.Range("I7").Value = 0.1592492032139157
$ws.Range("J7").Value = 0.1592492032139158
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01994566666666667
$ws.Range("N7").Value = 0.059837
$ws.Range("O7").Value = 0.01079376578266031
$ws.Range("P7").Value = 0.01079376578266031
$ws.Range("Q7").Value = 0.3910511238524444
$ws.Range("R7").Value = 3.519460114672
$ws.Range("S7").Value = 0.001718898600566282
$ws.Range("T7").Value = 0.001718898600566282

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Trf"
$ws.Range("C8").Value = "Tfr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.570446666666667
$ws.Range("H8").Value = 4.71134
$ws.Range("I8").Value = 0.01275602843261105
$ws.Range("J8").Value = 0.01275602843261106
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.771067666666666
$ws.Range("N8").Value = 5.313203
$ws.Range("O8").Value = 0.9584282089297276
$ws.Range("P8").Value = 0.9584282089297276
$ws.Range("Q8").Value = 2.781367313557777
$ws.Range("R8").Value = 25.03230582202
$ws.Range("S8").Value = 0.01222573748372409
$ws.Range("T8").Value = 0.01222573748372409

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Trf"
$ws.Range("C9").Value = "Tfr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.570446666666667
$ws.Range("H9").Value = 4.71134
$ws.Range("I9").Value = 0.01275602843261105
$ws.Range("J9").Value = 0.01275602843261106
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05687433333333333
$ws.Range("N9").Value = 0.170623
$ws.Range("O9").Value = 0.03077802528761218
$ws.Range("P9").Value = 0.03077802528761218
$ws.Range("Q9").Value = 0.08931810720222222
$ws.Range("R9").Value = 0.80386296482
$ws.Range("S9").Value = 0.000392605365668403
$ws.Range("T9").Value = 0.000392605365668403

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Trf"
$ws.Range("C10").Value = "Tfr2"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.570446666666667
$ws.Range("H10").Value = 4.71134
$ws.Range("I10").Value = 0.01275602843261105
$ws.Range("J10").Value = 0.01275602843261106
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.01994566666666667
$ws.Range("N10").Value = 0.059837
$ws.Range("O10").Value = 0.01079376578266031
$ws.Range("P10").Value = 0.01079376578266031
$ws.Range("Q10").Value = 0.03132360573111111
$ws.Range("R10").Value = 0.28191245158
$ws.Range("S10").Value = 0.0001376855832185592
$ws.Range("T10").Value = 0.0001376855832185592

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Trf"
$ws.Range("C11").Value = "Tfr2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 101.3157653333333
$ws.Range("H11").Value = 303.947296
$ws.Range("I11").Value = 0.8229421671522854
$ws.Range("J11").Value = 0.8229421671522854
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.771067666666666
$ws.Range("N11").Value = 5.313203
$ws.Range("O11").Value = 0.9584282089297276
$ws.Range("P11").Value = 0.9584282089297276
$ws.Range("Q11").Value = 179.4370761054542
$ws.Range("R11").Value = 1614.933684949088
$ws.Range("S11").Value = 0.7887309873165134
$ws.Range("T11").Value = 0.7887309873165134

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Trf"
$ws.Range("C12").Value = "Tfr2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 101.3157653333333
$ws.Range("H12").Value = 303.947296
$ws.Range("I12").Value = 0.8229421671522854
$ws.Range("J12").Value = 0.8229421671522854
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.05687433333333333
$ws.Range("N12").Value = 0.170623
$ws.Range("O12").Value = 0.03077802528761218
$ws.Range("P12").Value = 0.03077802528761218
$ws.Range("Q12").Value = 5.762266609489777
$ws.Range("R12").Value = 51.860399485408
$ws.Range("S12").Value = 0.02532853483085541
$ws.Range("T12").Value = 0.02532853483085541

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Trf"
$ws.Range("C13").Value = "Tfr2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 101.3157653333333
$ws.Range("H13").Value = 303.947296
$ws.Range("I13").Value = 0.8229421671522854
$ws.Range("J13").Value = 0.8229421671522854
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.01994566666666667
$ws.Range("N13").Value = 0.059837
$ws.Range("O13").Value = 0.01079376578266031
$ws.Range("P13").Value = 0.01079376578266031
$ws.Range("Q13").Value = 2.020810483416889
$ws.Range("R13").Value = 18.187294350752
$ws.Range("S13").Value = 0.00888264500491666
$ws.Range("T13").Value = 0.008882645004916659
